$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "top flight xl womens",
    "top generators",
    "top hike sit up",
    "top leggins",
    "top of hip pain",
    "top tights for women",
    "tops for women long length",
    "tops to go over leggings",
    "tops to go with leggings",
    "torn leggings",
    "torn leggings for women",
    "torn muscle",
    "torn pants for women",
    "torn tights for women",
    "torso target",
    "torso track 2 exercise machine",
    "torso track exercise machine",
    "total 90 soccer",
    "total body workout chart",
    "total control training ball",
    "total joint",
    "total support hose",
    "touch compression",
    "tough athletics leggings",
    "tough pants",
    "tough yoga ball",
    "town dump game",
    "track and field clothes for women",
    "track and field compression",
    "track and field gear",
    "track and field pants",
    "track and field shorts",
    "track and field tights for men",
    "track bike frame",
    "track capri for women",
    "track compression shorts men",
    "track compression shorts women",
    "track flights",
    "track frame",
    "track gear for women",
    "track leggings",
    "track leggings men",
    "track leggings women",
    "track runner",
    "track running tights",
    "track shorts",
    "track skirt",
    "track tights",
    "track tights for boys",
    "track tights for men",
    "track tights for women",
    "track tights men",
    "track training equipment",
    "track travel",
    "tracking package from amazon",
    "tracking pants for men",
    "trail pants women",
    "trail runners women",
    "trail running light",
    "trail running pants",
    "trail running pants men",
    "trail running tights",
    "trail tights",
    "training athletic pants",
    "training basketball heavy",
    "training capris women",
    "training clothes",
    "training compression",
    "training compression men",
    "training exercises for home b-80",
    "training gear for women",
    "training leggings",
    "training leggings for men",
    "training leggings for women",
    "training pant men",
    "training pants for man",
    "training pants for women",
    "training pants women",
    "training pennies",
    "training running",
    "training shells",
    "training shorts for women",
    "training support",
    "training tights",
    "training tights men",
    "training tights women",
    "training to be all might",
    "training waist",
    "training waist for women",
    "training wear",
    "training wear for women",
    "training women clothes",
    "training women pants",
    "trainjng pants",
    "trainning leggings for women",
    "travel capri pants women",
    "travel capris",
    "travel clothes for women over 50",
    "travel compression tights",
    "travel exercise"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
